$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "38.341.87"
$ws.Range("E2").Value = "  +1.67%  "

Set-TextValue $ws.Range("D3") "2.095.83"
$ws.Range("E3").Value = "  +3.47%  "

$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue $ws.Range("D5") "228.48"

$ws.Range("E6").Value = "  +1.35%  "

Set-TextValue $ws.Range("D7") "61.11"
$ws.Range("E7").Value = "  +1.90%  "

$ws.Range("E8").Value = "  -0.01%  "

Set-TextValue $ws.Range("D9") "0.380"
$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("E10").Value = "  +2.84%  "

$ws.Range("E11").Value = "  +0.39%  "

Set-TextValue $ws.Range("D12") "2.405.72"
$ws.Range("E12").Value = "  +3.43%  "

Set-TextValue $ws.Range("D13") "14.77"
$ws.Range("E13").Value = "  +2.84%  "

Set-TextValue $ws.Range("D14") "22.34"
$ws.Range("E14").Value = "  +6.43%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D15") "0.776"
$ws.Range("E15").Value = "  +2.51%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D16") "5.44"
$ws.Range("E16").Value = "  +5.45%  "

Set-TextValue $ws.Range("D17") "2.099.56"
$ws.Range("E17").Value = "  +3.75%  "

Set-TextValue $ws.Range("D18") "38.297.80"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("E19").Value = "  +1.38%  "

Set-TextValue $ws.Range("D20") "6.01"
$ws.Range("E20").Value = "  +1.90%  "

Set-TextValue $ws.Range("D21") "0.0₃0833"
$ws.Range("E21").Value = "  +1.41%  "

Set-TextValue $ws.Range("D22") "225.18"
$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("E23").Value = "  -0.05%  "

Set-TextValue $ws.Range("D24") "2.40"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("E25").Value = "  +2.81%  "

Set-TextValue $ws.Range("D26") "169.81"
$ws.Range("E26").Value = "  +1.46%  "

Set-TextValue $ws.Range("D27") "9.42"
$ws.Range("E27").Value = "  +1.86%  "

Set-TextValue $ws.Range("D28") "0.130"
$ws.Range("E28").Value = "  +1.05%  "

Set-TextValue $ws.Range("D29") "19.00"
$ws.Range("E29").Value = "  +1.25%  "

$ws.Range("E30").Value = "  +9.03%  "

Set-TextValue $ws.Range("D31") "0.119"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  +6.96%  "

Set-TextValue $ws.Range("D33") "4.77"
$ws.Range("E33").Value = "  +7.00%  "

Set-TextValue $ws.Range("D34") "4.44"
$ws.Range("E34").Value = "  +1.48%  "

$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("E36").Value = "  +1.09%  "

$ws.Range("E37").Value = "  +3.76%  "

$ws.Range("E38").Value = "  +4.00%  "

$ws.Range("E39").Value = "  +0.03%  "

Set-TextValue $ws.Range("D40") "18.21"
$ws.Range("E40").Value = "  +1.63%  "

Set-TextValue $ws.Range("D41") "1.541.11"
$ws.Range("E41").Value = "  +0.45%  "

Set-TextValue $ws.Range("D42") "99.77"
$ws.Range("E42").Value = "  +4.57%  "

$ws.Range("E43").Value = "  +1.80%  "

$ws.Range("E44").Value = "  +0.99%  "

Set-TextValue $ws.Range("D45") "0.0908"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("E46").Value = "  +1.53%  "

$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("E48").Value = "  +5.99%  "

$ws.Range("E49").Value = "  +3.00%  "

$ws.Range("E50").Value = "  +0.80%  "

Set-TextValue $ws.Range("D51") "2.290.93"
$ws.Range("E51").Value = "  +3.38%  "
